# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for rows 2-38 with the newly
# recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 5
    3  = 2
    4  = 4
    5  = 1
    6  = 1
    7  = 4
    8  = 5
    9  = 5
    10 = 5
    11 = 4
    12 = 9
    13 = 6
    14 = 8
    15 = 7
    16 = 3
    17 = 4
    18 = 8
    19 = 0
    20 = 6
    21 = 4
    22 = 2
    23 = 5
    24 = 2
    25 = 5
    26 = 5
    27 = 2
    28 = 5
    29 = 4
    30 = 3
    31 = 6
    32 = 3
    33 = 2
    34 = 5
    35 = 2
    36 = 7
    37 = 3
    38 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
